# Correction: scale Avg_Throughput_Before (col I) and Avg_Throughput_After (col J)
# for data rows 2-22 by a factor of 10 (unit fix from Kbps -> bps, per commit "Correction").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

for ($r = 2; $r -le 22; $r++) {
    $ws.Cells.Item($r, 9).Value2  = $ws.Cells.Item($r, 9).Value2  * 10   # column I
    $ws.Cells.Item($r, 10).Value2 = $ws.Cells.Item($r, 10).Value2 * 10  # column J
}

# Widen the (currently empty) helper columns L:M to fit content, matching the
# author's follow-up formatting tweak (best-fit column width -> stored width 10).
$ws.Range("L1:M1").EntireColumn.ColumnWidth = 9.14

# Move/save the active selection to N20, matching the author's last cursor position.
$ws.Range("N20").Select()
